# "Generate Report for Archive"
# - Status text "Ready for handoff" -> "In Translation" (Overview!E2:F2, zh-cn!C2, de-de!C2)
# - Narrower Status/zh-cn/de-de columns on Overview (E:F) and the Status column (C)
#   on the zh-cn / de-de detail sheets (width ~17.22 chars -> ~13.41 chars)

$wb = $excel.ActiveWorkbook

# Update the "Status" values wherever they appear, in every sheet, using
# Find & Replace so the shared-string content is edited in place rather than
# each cell being rewritten with a brand-new literal.
foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Narrow the Status columns.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
